$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("comforter-cda")

# Fill in the missing Start Time / End Time for the existing rows 89 and 90
$ws.Range("B89").Value = 0.80555555555555547
$ws.Range("C89").Value = 0.99930555555555556

$ws.Range("B90").Value = 0
$ws.Range("C90").Value = 0.36458333333333331

# Add a new daily record row (row 91): date + the three calculated columns.
# (Start Time / End Time are left blank for this row, matching the source data.)
# Copy the number formats from row 90 cell-by-cell so no stray cells get
# created in columns B/C (which must stay empty on row 91).
$ws.Range("A90").Copy() | Out-Null
$ws.Range("A91").PasteSpecial(-4122) | Out-Null

$ws.Range("D90").Copy() | Out-Null
$ws.Range("D91").PasteSpecial(-4122) | Out-Null

$ws.Range("E90").Copy() | Out-Null
$ws.Range("E91").PasteSpecial(-4122) | Out-Null

$ws.Range("F90").Copy() | Out-Null
$ws.Range("F91").PasteSpecial(-4122) | Out-Null

$ws.Range("A91").Value = 43415
$ws.Range("D91").Formula = "=(C91-B91)* 1440"
$ws.Range("E91").Formula = "=IF(C91>B91, (C91-B91)*1440, (B91-C91)*1440)"
$ws.Range("F91").Formula = "=ABS((C91-B91)*1440)"

# Resize the table to include the new row
$table = $ws.ListObjects.Item("comforter_cda_table")
$table.Resize($ws.Range("A1:F91"))

$ws.Range("B91").Select() | Out-Null

$excel.CutCopyMode = $false
